$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")

# Update the Version and Date values in the metadata table.
$ws1.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$ws1.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new "Jurisdiction" property row right after "Contact" (row 10),
# pushing Description/Purpose/Copyright/etc. down by one row.
$ws1.Rows.Item(11).Insert()

# Copy formatting from the row that is now below (the old row 11, "Description",
# now at row 12) so the new row matches the existing table styling exactly.
$ws1.Range("A12:B12").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)

$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = ""
